# Updated cryptos list on Sat Mar 30 22:28:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay stored as text
# (matching the source feed formatting, e.g. trailing zeros/thousand dots).
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"

# Apply the updated price / volume(1h) figures scraped this run.
$ws.Range("D2").Value = "69.635.36"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.512.31"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "603.62"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "196.01"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.201"
$ws.Range("E9").Value = "  -5.48%  "
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "53.20"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "0.0000297"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "9.44"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "4.057.83"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "597.45"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "69.751.85"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "19.01"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "12.65"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "3.507.19"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "0.984"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "18.04"
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Value = "101.81"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").Value = "4.63"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").Value = "3.12"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").Value = "10.76"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "9.48"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "33.00"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "4.27"
$ws.Range("E30").Value = "  +8.67%  "
$ws.Range("D31").Value = "6.99"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "12.31"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "63.14"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "3.744.39"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "0.0₃0810"
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "3.64"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "0.389"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "36.14"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "491.42"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").Value = "0.0450"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "8.39"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("D50").Value = "0.000243"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "1.33"
$ws.Range("E51").Value = "  -1.65%  "
